$d = $word.ActiveDocument

# --- Edit 1: "Outcome" paragraph -------------------------------------
# Old: "Predicted accuracy was 76% and we need extend offers for other users."
# New: "Predicted accuracy was 76% and we need not extend offer to other
#       users and we can target specific users by this anaysis."
$found1 = $d.Content.Find.Execute(
    "Predicted accuracy was 76% and we need extend offers for other users.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Predicted accuracy was 76% and we need not extend offer to other users and we can target specific users by this anaysis.",
    2)

# --- Edit 2: "Tools Used" / Python paragraph --------------------------
# Collapse the "Python – Numpy and Pandas..." run (currently split across
# three runs around a spell-checker proofErr tag for "Numpy") into a
# single clean run with the same wording.
$found2 = $d.Content.Find.Execute(
    "Python – Numpy and Pandas - Data Manipulation and formatting",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Python – Numpy and Pandas - Data Manipulation and formatting",
    2)

Write-Output "Edit1: $found1, Edit2: $found2"
Write-Output $d.Content.Text
